# Auto-generated Excel COM-interop script applying scheduled market-data refresh
# to the per-Leve profit columns (H..N) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 181.875
$ws.Range("I28").Value = 209.16667
$ws.Range("K28").Value = 209.16667
$ws.Range("M28").Value = 275.83333
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H53").Value = 948.46155
$ws.Range("I53").Value = 686
$ws.Range("J53").Value = 2050.8
$ws.Range("K53").Value = 686
$ws.Range("L53").Value = 2050.8
$ws.Range("M53").Value = -49
$ws.Range("N53").Value = -3324.8
$ws.Range("H80").Value = 1371.5555
$ws.Range("J80").Value = 963.4286
$ws.Range("L80").Value = 2890.2858
$ws.Range("N80").Value = -4886.2858
$ws.Range("H83").Value = 1371.5555
$ws.Range("J83").Value = 963.4286
$ws.Range("L83").Value = 8670.857399999999
$ws.Range("N83").Value = -18654.8574
$ws.Range("H88").Value = 1374629.6
$ws.Range("I88").Value = 999.5
$ws.Range("J88").Value = 1767095.4
$ws.Range("K88").Value = 999.5
$ws.Range("L88").Value = 1767095.4
$ws.Range("M88").Value = -593.5
$ws.Range("N88").Value = -1767907.4
$ws.Range("H91").Value = 1374629.6
$ws.Range("I91").Value = 999.5
$ws.Range("J91").Value = 1767095.4
$ws.Range("K91").Value = 999.5
$ws.Range("L91").Value = 1767095.4
$ws.Range("M91").Value = 404.5
$ws.Range("N91").Value = -1769903.4
$ws.Range("H112").Value = 2808.6667
$ws.Range("J112").Value = 2986.2666
$ws.Range("L112").Value = 8958.799800000001
$ws.Range("N112").Value = -11174.7998
$ws.Range("H116").Value = 3947.8235
$ws.Range("I116").Value = 3509.7273
$ws.Range("K116").Value = 3509.7273
$ws.Range("M116").Value = -67.72730000000001
$ws.Range("H129").Value = 869.0227
$ws.Range("J129").Value = 899.5122
$ws.Range("L129").Value = 2698.5366
$ws.Range("N129").Value = -12698.5366
$ws.Range("H137").Value = 1368.3334
$ws.Range("I137").Value = 1160.6666
$ws.Range("J137").Value = 1576
$ws.Range("K137").Value = 3481.9998
$ws.Range("L137").Value = 4728
$ws.Range("M137").Value = -931.9998000000001
$ws.Range("N137").Value = -9828
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3499.76
$ws.Range("I32").Value = 3158.9546
$ws.Range("J32").Value = 5999
$ws.Range("K32").Value = 3158.9546
$ws.Range("L32").Value = 5999
$ws.Range("M32").Value = -2871.9546
$ws.Range("N32").Value = -6573
$ws.Range("H37").Value = 11038
$ws.Range("J37").Value = 11038
$ws.Range("L37").Value = 11038
$ws.Range("N37").Value = -11584
$ws.Range("H61").Value = 2223.375
$ws.Range("I61").Value = 1697
$ws.Range("J61").Value = 2749.75
$ws.Range("K61").Value = 1697
$ws.Range("L61").Value = 2749.75
$ws.Range("M61").Value = -1485
$ws.Range("N61").Value = -3173.75
$ws.Range("H74").Value = 1664.2222
$ws.Range("I74").Value = 1006
$ws.Range("K74").Value = 1006
$ws.Range("M74").Value = -132
$ws.Range("H77").Value = 1664.2222
$ws.Range("I77").Value = 1006
$ws.Range("K77").Value = 5030
$ws.Range("M77").Value = -662
$ws.Range("H88").Value = 2588.4
$ws.Range("I88").Value = 2238.75
$ws.Range("J88").Value = 2675.8125
$ws.Range("K88").Value = 2238.75
$ws.Range("L88").Value = 2675.8125
$ws.Range("M88").Value = -1832.75
$ws.Range("N88").Value = -3487.8125
$ws.Range("H91").Value = 2588.4
$ws.Range("I91").Value = 2238.75
$ws.Range("J91").Value = 2675.8125
$ws.Range("K91").Value = 2238.75
$ws.Range("L91").Value = 2675.8125
$ws.Range("M91").Value = -834.75
$ws.Range("N91").Value = -5483.8125
$ws.Range("H136").Value = 2223.375
$ws.Range("I136").Value = 1697
$ws.Range("J136").Value = 2749.75
$ws.Range("K136").Value = 5091
$ws.Range("L136").Value = 8249.25
$ws.Range("M136").Value = -2541
$ws.Range("N136").Value = -13349.25
$ws.Range("H139").Value = 33190.625
$ws.Range("J139").Value = 33190.625
$ws.Range("L139").Value = 33190.625
$ws.Range("N139").Value = -43470.625
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3678.577
$ws.Range("J86").Value = 2893.7144
$ws.Range("L86").Value = 2893.7144
$ws.Range("N86").Value = -5139.7144
$ws.Range("H89").Value = 3678.577
$ws.Range("J89").Value = 2893.7144
$ws.Range("L89").Value = 14468.572
$ws.Range("N89").Value = -25700.572
$ws.Range("H134").Value = 8208.579
$ws.Range("I134").Value = 1163.6666
$ws.Range("K134").Value = 3490.9998
$ws.Range("M134").Value = -955.9998000000001
$ws.Range("H138").Value = 32879
$ws.Range("J138").Value = 32879
$ws.Range("L138").Value = 32879
$ws.Range("N138").Value = -43159
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1468.6
$ws.Range("I31").Value = 1335.4348
$ws.Range("K31").Value = 1335.4348
$ws.Range("M31").Value = -1040.4348
$ws.Range("H34").Value = 1468.6
$ws.Range("I34").Value = 1335.4348
$ws.Range("K34").Value = 1335.4348
$ws.Range("M34").Value = -1133.4348
$ws.Range("H122").Value = 1022.2273
$ws.Range("I122").Value = 896.25
$ws.Range("J122").Value = 1358.1666
$ws.Range("K122").Value = 2688.75
$ws.Range("L122").Value = 4074.4998
$ws.Range("M122").Value = -238.75
$ws.Range("N122").Value = -8974.4998
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 36.307693
$ws.Range("I2").Value = 12.5
$ws.Range("K2").Value = 75
$ws.Range("M2").Value = 38
$ws.Range("H62").Value = 3157
$ws.Range("I62").Value = 300
$ws.Range("J62").Value = 6014
$ws.Range("K62").Value = 900
$ws.Range("L62").Value = 18042
$ws.Range("M62").Value = -214
$ws.Range("N62").Value = -19414
$ws.Range("H65").Value = 3157
$ws.Range("I65").Value = 300
$ws.Range("J65").Value = 6014
$ws.Range("K65").Value = 2700
$ws.Range("L65").Value = 54126
$ws.Range("M65").Value = 732
$ws.Range("N65").Value = -60990
$ws.Range("H68").Value = 951.125
$ws.Range("I68").Value = 577.5
$ws.Range("J68").Value = 1175.3
$ws.Range("K68").Value = 1732.5
$ws.Range("L68").Value = 3525.9
$ws.Range("M68").Value = -921.5
$ws.Range("N68").Value = -5147.9
$ws.Range("H71").Value = 951.125
$ws.Range("I71").Value = 577.5
$ws.Range("J71").Value = 1175.3
$ws.Range("K71").Value = 5197.5
$ws.Range("L71").Value = 10577.7
$ws.Range("M71").Value = -1141.5
$ws.Range("N71").Value = -18689.7
$ws.Range("H132").Value = 919
$ws.Range("J132").Value = 1038.75
$ws.Range("L132").Value = 9348.75
$ws.Range("N132").Value = -14408.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3731.5
$ws.Range("I132").Value = 4047.25
$ws.Range("K132").Value = 12141.75
$ws.Range("M132").Value = -9611.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 18732.55
$ws.Range("I132").Value = 1113.5
$ws.Range("J132").Value = 47563.727
$ws.Range("K132").Value = 3340.5
$ws.Range("L132").Value = 142691.181
$ws.Range("M132").Value = -810.5
$ws.Range("N132").Value = -147751.181
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").ClearContents()
$ws.Range("H112").Value = 28600
$ws.Range("J112").Value = 28600
$ws.Range("L112").Value = 28600
$ws.Range("N112").Value = -31554
$ws.Range("H113").Value = 338.38095
$ws.Range("I113").Value = 207.05882
$ws.Range("K113").Value = 621.17646
$ws.Range("M113").Value = 1548.82354
